# The diff shows a new weekly price record being inserted before the
# existing row 47 (pushing rows 47-69 down to 48-70), with the sheet's
# used-range dimension growing from A1:R69 to A1:R70.
#
# Insert a blank row at position 47 (shifts rows 47-69 -> 48-70, carrying
# their formatting/styles down with them), then populate the new row 47
# with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(47).Insert()

$ws.Cells.Item(47, 1).Value = 4
$ws.Cells.Item(47, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(47, 3).Value = "Los Lagos"
$ws.Cells.Item(47, 4).Value = 44488
$ws.Cells.Item(47, 5).Value = 10
$ws.Cells.Item(47, 6).Value = 100112052
$ws.Cells.Item(47, 7).Value = "Albahaca"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 90
$ws.Cells.Item(47, 11).Value = 6000
$ws.Cells.Item(47, 12).Value = 6000
$ws.Cells.Item(47, 13).Value = 6000
$ws.Cells.Item(47, 14).Value = "`$/paquete"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 6000
$ws.Cells.Item(47, 17).Value = 1
$ws.Cells.Item(47, 18).Value = "Hortaliza"
